$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

# Copy the header style (from H1, which already has the bold/border/centered style) to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
